$d = $word.ActiveDocument

# Paragraph 1 (Title) currently starts with a run containing just a space,
# followed by the "Algoritmo para trocar um pneu de um carro" run.
# We need to insert a new run "03-" (Arial, color 3C4043, spacing 3) at the
# very start of the paragraph, followed immediately by the "_GoBack" bookmark
# (which today lives, empty, in the following paragraph).

$p1 = $d.Paragraphs.Item(1)
$insPoint = $p1.Range.Duplicate
$insPoint.Collapse(1)          # wdCollapseStart
$insPoint.InsertBefore("03-")

# Range covering exactly the newly inserted "03-" text (3 characters at the
# very start of the document/paragraph).
$newRun = $d.Range(0, 3)
$newRun.Font.NameAscii = "Arial"
$newRun.Font.NameOther = "Arial"
$newRun.Font.NameBi = "Arial"
$newRun.Font.Color = 4407356    # RGB(3C,40,43) packed as BGR for wdColor
$newRun.Font.Spacing = 0.15     # 3 (twentieths of a point) == 0.15 pt

# Re-home the "_GoBack" bookmark right after the new "03-" run (Word only
# allows one bookmark with a given name, so adding it here removes the
# existing, empty one that previously lived alone in paragraph 2).
$bmPoint = $d.Range(3, 3)
$d.Bookmarks.Add("_GoBack", $bmPoint)
